# Update the "loading_percent" results table (rows 2-25, cols B..N) with
# the recomputed loading percentages for the 380 kV case.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 18.09656102687515
$ws.Range("C2").Value = 10.45841578630965
$ws.Range("D2").Value = 5.044523315965096
$ws.Range("F2").Value = 28.23900646101855
$ws.Range("G2").Value = 36.21740050928798
$ws.Range("H2").Value = 15.63768947374823
$ws.Range("L2").Value = 10.70357779574878
$ws.Range("N2").Value = 17.71092062368574

$ws.Range("B3").Value = 17.55081334527085
$ws.Range("C3").Value = 10.20244018724204
$ws.Range("D3").Value = 5.057315094279921
$ws.Range("F3").Value = 28.07737362686862
$ws.Range("G3").Value = 35.86057660481368
$ws.Range("H3").Value = 15.65033159118782
$ws.Range("L3").Value = 10.6789024190537
$ws.Range("N3").Value = 17.78362784503239

$ws.Range("B4").Value = 17.21131276488243
$ws.Range("C4").Value = 10.04032588521472
$ws.Range("D4").Value = 5.065512849162102
$ws.Range("F4").Value = 27.98748956236277
$ws.Range("G4").Value = 35.65452867888359
$ws.Range("H4").Value = 15.66202622670902
$ws.Range("L4").Value = 10.66602619474273
$ws.Range("N4").Value = 17.83022434970528

$ws.Range("B5").Value = 17.07208475089453
$ws.Range("C5").Value = 9.973074127493952
$ws.Range("D5").Value = 5.068940258457347
$ws.Range("F5").Value = 27.95324336417967
$ws.Range("G5").Value = 35.57393481202429
$ws.Range("H5").Value = 15.66777808142574
$ws.Range("L5").Value = 10.66135422404347
$ws.Range("N5").Value = 17.84970593801406

$ws.Range("B6").Value = 17.04891997420786
$ws.Range("C6").Value = 9.961836882980602
$ws.Range("D6").Value = 5.069514627981034
$ws.Range("F6").Value = 27.94770144851478
$ws.Range("G6").Value = 35.56075847470863
$ws.Range("H6").Value = 15.66879265069781
$ws.Range("L6").Value = 10.66061327148215
$ws.Range("N6").Value = 17.85297067759875

$ws.Range("B7").Value = 17.20943833076921
$ws.Range("C7").Value = 10.03942364350024
$ws.Range("D7").Value = 5.065558720652985
$ws.Range("F7").Value = 27.98701802568064
$ws.Range("G7").Value = 35.65342799171523
$ws.Range("H7").Value = 15.66209980903304
$ws.Range("L7").Value = 10.66596085408477
$ws.Range("N7").Value = 17.83048508619882

$ws.Range("B8").Value = 17.90943576993095
$ws.Range("C8").Value = 10.37121305516952
$ws.Range("D8").Value = 5.048862839028669
$ws.Range("F8").Value = 28.18135163447363
$ws.Range("G8").Value = 36.0917185508642
$ws.Range("H8").Value = 15.64123088131725
$ws.Range("L8").Value = 10.69459974108618
$ws.Range("N8").Value = 17.73558560927927

$ws.Range("B9").Value = 19.23781043110991
$ws.Range("C9").Value = 10.98042433022799
$ws.Range("D9").Value = 5.018831343927556
$ws.Range("F9").Value = 28.63522102173927
$ws.Range("G9").Value = 37.05023065521826
$ws.Range("H9").Value = 15.63160635566445
$ws.Range("L9").Value = 10.7686432424215
$ws.Range("N9").Value = 17.56490881563089

$ws.Range("B10").Value = 20.17538767169534
$ws.Range("C10").Value = 11.40008626163633
$ws.Range("D10").Value = 4.998395125318728
$ws.Range("F10").Value = 29.01078224462694
$ws.Range("G10").Value = 37.80833839903385
$ws.Range("H10").Value = 15.64372220790186
$ws.Range("L10").Value = 10.83369172491665
$ws.Range("N10").Value = 17.44879509800196

$ws.Range("B11").Value = 20.59150889117659
$ws.Range("C11").Value = 11.58444690228107
$ws.Range("D11").Value = 4.989446658366421
$ws.Range("F11").Value = 29.19021610685331
$ws.Range("G11").Value = 38.16343118507812
$ws.Range("H11").Value = 15.65341247299192
$ws.Range("L11").Value = 10.8655317140496
$ws.Range("N11").Value = 17.39796218425152

$ws.Range("B12").Value = 20.74743863256111
$ws.Range("C12").Value = 11.65328187392522
$ws.Range("D12").Value = 4.986107786337311
$ws.Range("F12").Value = 29.25934615799633
$ws.Range("G12").Value = 38.29923589749207
$ws.Range("H12").Value = 15.65768285910966
$ws.Range("L12").Value = 10.87790584637968
$ws.Range("N12").Value = 17.37899702572038

$ws.Range("B13").Value = 20.71393185929248
$ws.Range("C13").Value = 11.63850110097509
$ws.Range("D13").Value = 4.98682466619949
$ws.Range("F13").Value = 29.24440604876942
$ws.Range("G13").Value = 38.26993049790478
$ws.Range("H13").Value = 15.65673643604444
$ws.Range("L13").Value = 10.87522686291598
$ws.Range("N13").Value = 17.3830689003814

$ws.Range("B14").Value = 20.60437113885803
$ws.Range("C14").Value = 11.59012980605331
$ws.Range("D14").Value = 4.989170972825767
$ws.Range("F14").Value = 29.19588009386323
$ws.Range("G14").Value = 38.17457785077521
$ws.Range("H14").Value = 15.65375175940624
$ws.Range("L14").Value = 10.86654342613747
$ws.Range("N14").Value = 17.39639622379755

$ws.Range("B15").Value = 20.5370433686092
$ws.Range("C15").Value = 11.56037253053252
$ws.Range("D15").Value = 4.990614618327275
$ws.Range("F15").Value = 29.16630890797081
$ws.Range("G15").Value = 38.11634200650278
$ws.Range("H15").Value = 15.65200180334922
$ws.Range("L15").Value = 10.86126565775034
$ws.Range("N15").Value = 17.40459655085308

$ws.Range("B16").Value = 20.14797074491847
$ws.Range("C16").Value = 11.3879028909346
$ws.Range("D16").Value = 4.998986904094313
$ws.Range("F16").Value = 28.99922420902887
$ws.Range("G16").Value = 37.78532681458391
$ws.Range("H16").Value = 15.6431730251634
$ws.Range("L16").Value = 10.83165562828297
$ws.Range("N16").Value = 17.45215700420837

$ws.Range("B17").Value = 19.90651323573619
$ws.Range("C17").Value = 11.280394325211
$ws.Range("D17").Value = 5.004211940381687
$ws.Range("F17").Value = 28.89888641230293
$ws.Range("G17").Value = 37.58478585256096
$ws.Range("H17").Value = 15.63882729766255
$ws.Range("L17").Value = 10.81406248236135
$ws.Range("N17").Value = 17.48184172629237

$ws.Range("B18").Value = 19.76666366080921
$ws.Range("C18").Value = 11.21794412823033
$ws.Range("D18").Value = 5.007250024852923
$ws.Range("F18").Value = 28.8419865707055
$ws.Range("G18").Value = 37.47041055810389
$ws.Range("H18").Value = 15.63672110454115
$ws.Range("L18").Value = 10.80415548065922
$ws.Range("N18").Value = 17.49910279428152

$ws.Range("B19").Value = 19.71915137595466
$ws.Range("C19").Value = 11.19669527823256
$ws.Range("D19").Value = 5.008284308607176
$ws.Range("F19").Value = 28.82286225173864
$ws.Range("G19").Value = 37.43185571763401
$ws.Range("H19").Value = 15.63607553360323
$ws.Range("L19").Value = 10.80083775750812
$ws.Range("N19").Value = 17.50497929857809

$ws.Range("B20").Value = 19.93231830532058
$ws.Range("C20").Value = 11.29190263382069
$ws.Range("D20").Value = 5.003652336029591
$ws.Range("F20").Value = 28.9094838969209
$ws.Range("G20").Value = 37.60603427785438
$ws.Range("H20").Value = 15.63924919546512
$ws.Range("L20").Value = 10.81591339502548
$ws.Range("N20").Value = 17.47866237465131

$ws.Range("B21").Value = 20.63659759863369
$ws.Range("C21").Value = 11.60436446051654
$ws.Range("D21").Value = 4.988480458992785
$ws.Range("F21").Value = 29.21010168879615
$ws.Range("G21").Value = 38.20254995993311
$ws.Range("H21").Value = 15.65461212595905
$ws.Range("L21").Value = 10.86908541150969
$ws.Range("N21").Value = 17.39247396710303

$ws.Range("B22").Value = 21.08722874956784
$ws.Range("C22").Value = 11.80285860143316
$ws.Range("D22").Value = 4.978854404527437
$ws.Range("F22").Value = 29.41343806702017
$ws.Range("G22").Value = 38.60014861315124
$ws.Range("H22").Value = 15.6681547619358
$ws.Range("L22").Value = 10.90568093725988
$ws.Range("N22").Value = 17.33780040448949

$ws.Range("B23").Value = 20.84764840009717
$ws.Range("C23").Value = 11.69745310984417
$ws.Range("D23").Value = 4.983965617601601
$ws.Range("F23").Value = 29.30430356273917
$ws.Range("G23").Value = 38.3872783830096
$ws.Range("H23").Value = 15.66060650559635
$ws.Range("L23").Value = 10.88598267464678
$ws.Range("N23").Value = 17.36682978438887

$ws.Range("B24").Value = 19.92065503740006
$ws.Range("C24").Value = 11.28670172466003
$ws.Range("D24").Value = 5.003905226769597
$ws.Range("F24").Value = 28.9046903227805
$ws.Range("G24").Value = 37.59642499415403
$ws.Range("H24").Value = 15.639057233716
$ws.Range("L24").Value = 10.81507595050589
$ws.Range("N24").Value = 17.48009915220515

$ws.Range("B25").Value = 18.88443148735831
$ws.Range("C25").Value = 10.82034181557001
$ws.Range("D25").Value = 5.026668092078125
$ws.Range("F25").Value = 28.50486411042888
$ws.Range("G25").Value = 36.78097190414491
$ws.Range("H25").Value = 15.63084625317993
$ws.Range("L25").Value = 10.74672209805279
$ws.Range("N25").Value = 17.60944278235863
